$d = $word.ActiveDocument

# The "Requisitos" section ends with the line:
#   "LOQ4233: Gestão de Negócios (Requisito fraco)"
# Immediately after it the page footer (an empty spacer paragraph, the
# "Ver no Jupiter..." line, and the "© 2020 ..." copyright line) is being
# dropped from this rebuilt page. Locate that anchor paragraph and remove
# the three paragraphs that follow it, leaving the trailing blank paragraph
# and the page-break paragraph untouched.

$count = $d.Paragraphs.Count
$anchorIndex = -1
for ($i = 1; $i -le $count; $i++) {
    $t = $d.Paragraphs.Item($i).Range.Text
    if ($t -like "*LOQ4233*Requisito fraco*") {
        $anchorIndex = $i
        break
    }
}

if ($anchorIndex -lt 0) {
    throw "Could not locate the 'LOQ4233 ... (Requisito fraco)' anchor paragraph"
}

# Delete highest index first so earlier indices stay valid while we work.
$d.Paragraphs.Item($anchorIndex + 3).Range.Delete()
$d.Paragraphs.Item($anchorIndex + 2).Range.Delete()
$d.Paragraphs.Item($anchorIndex + 1).Range.Delete()
